$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the first sheet from "your data" to "Order list"
$ws.Name = "Order list"

# 2) Update the PickupDateTime column/header/cell number format from a
#    plain date ("[$-409]d-mmm-yy;@") to a date+time format
#    ("[$-409]m/d/yy h:mm AM/PM;@"). This affects the whole column R,
#    including the header (R1) and the sample data row (R2).
$ws.Columns("R:R").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
$ws.Range("R1:R2").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"

# 3) Update the view: scroll so column P is the first visible column and
#    select cell R2 (previously the view was scrolled to column S with R2
#    selected... now the pickup-date cell itself is selected).
$ws.Range("R2").Select()
